$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.918.33'
$ws.Range("E2").Value = '  -1.97%  '

# Row 3
$ws.Range("D3").Value = '1.899.78'
$ws.Range("E3").Value = '  -3.78%  '

# Row 4
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.31%  '

# Row 5
$ws.Range("D5").Value = '324.21'
$ws.Range("E5").Value = '  -0.83%  '

# Row 6
$ws.Range("E6").Value = '  -0.23%  '

# Row 7
$ws.Range("D7").Value = '0.4585'
$ws.Range("E7").Value = '  -1.62%  '

# Row 8
$ws.Range("D8").Value = '0.3816'
$ws.Range("E8").Value = '  -2.51%  '

# Row 9
$ws.Range("E9").Value = '  -2.98%  '

# Row 10
$ws.Range("D10").Value = '0.9759'
$ws.Range("E10").Value = '  -1.63%  '

# Row 11
$ws.Range("D11").Value = '22.04'
$ws.Range("E11").Value = '  -3.47%  '

# Row 12
$ws.Range("D12").Value = '1.894.63'
$ws.Range("E12").Value = '  -4.43%  '

# Row 13
$ws.Range("D13").Value = '6.935'
$ws.Range("E13").Value = '  -3.47%  '

# Row 14
$ws.Range("D14").Value = '5.641'
$ws.Range("E14").Value = '  -3.65%  '

# Row 15
$ws.Range("D15").Value = '0.07021'
$ws.Range("E15").Value = '  -0.79%  '

# Row 16
$ws.Range("E16").Value = '  -0.28%  '

# Row 17
$ws.Range("D17").Value = '83.55'
$ws.Range("E17").Value = '  -4.70%  '

# Row 18
$ws.Range("D18").Value = '0.000009456'
$ws.Range("E18").Value = '  -4.82%  '

# Row 19
$ws.Range("D19").Value = '16.60'
$ws.Range("E19").Value = '  -3.99%  '

# Row 21
$ws.Range("D21").Value = '28.903.63'
$ws.Range("E21").Value = '  -2.04%  '

# Row 22
$ws.Range("D22").Value = '5.292'
$ws.Range("E22").Value = '  -4.61%  '

# Row 23
$ws.Range("D23").Value = '10.85'
$ws.Range("E23").Value = '  -2.90%  '

# Row 24
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.141.28'
$ws.Range("E24").Value = '  -3.48%  '

# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.096'
$ws.Range("E25").Value = '  -0.67%  '

# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '157.99'
$ws.Range("E26").Value = '  -0.34%  '

# Row 27
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.97'
$ws.Range("E27").Value = '  -2.99%  '

# Row 28
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '5.633'
$ws.Range("E28").Value = '  -2.69%  '

# Row 29
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '117.45'
$ws.Range("E29").Value = '  -1.79%  '

# Row 30
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").Value = '1.834'
$ws.Range("E30").Value = '  -3.90%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.09247'
$ws.Range("E31").Value = '  -1.81%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = '0.8625'
$ws.Range("E32").Value = '  -3.23%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.081'
$ws.Range("E33").Value = '  -2.94%  '

# Row 34
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.240'
$ws.Range("E34").Value = '  -6.22%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '2.993'
$ws.Range("E35").Value = '  -6.15%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.05705'
$ws.Range("E36").Value = '  -2.02%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.147'
$ws.Range("E37").Value = '  -2.11%  '

# Row 38
$ws.Range("B38").Value = 'Frax'
$ws.Range("C38").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D38").Value = '1.002'
$ws.Range("E38").Value = '  -0.10%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02028'
$ws.Range("E39").Value = '  -3.58%  '

# Row 40
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5481'
$ws.Range("E40").Value = '  -4.19%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '7.385'
$ws.Range("E41").Value = '  -4.92%  '

# Row 42
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1750'
$ws.Range("E42").Value = '  -2.80%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '9.259'
$ws.Range("E43").Value = '  -4.05%  '

# Row 44
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '2.761'
$ws.Range("E44").Value = '  +0.01%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.5154'
$ws.Range("E45").Value = '  -3.79%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '11.28'
$ws.Range("E46").Value = '  -4.11%  '

# Row 47
$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").Value = '0.000002627'
$ws.Range("E47").Value = '  -17.15%  '

# Row 48
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.06807'
$ws.Range("E48").Value = '  -1.71%  '

# Row 49
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '2.069'
$ws.Range("E49").Value = '  -5.58%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '110.13'
$ws.Range("E50").Value = '  -3.50%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.770'
$ws.Range("E51").Value = '  -3.32%  '
